$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.734287
$ws.Range("H2").Value = 44.202861
$ws.Range("I2").Value = 0.4000023944294819
$ws.Range("J2").Value = 0.400002394429482
$ws.Range("M2").Value = 2.231113333333334
$ws.Range("N2").Value = 6.69334
$ws.Range("O2").Value = 0.01598125358798882
$ws.Range("P2").Value = 0.01598125358798882
$ws.Range("Q2").Value = 32.87386418286
$ws.Range("R2").Value = 295.86477764574
$ws.Range("S2").Value = 0.006392539701180276
$ws.Range("T2").Value = 0.006392539701180279

$ws.Range("G3").Value = 14.734287
$ws.Range("H3").Value = 44.202861
$ws.Range("I3").Value = 0.4000023944294819
$ws.Range("J3").Value = 0.400002394429482
$ws.Range("O3").Value = 0.1634493267640196
$ws.Range("P3").Value = 0.1634493267640195
$ws.Range("Q3").Value = 336.219617518533
$ws.Range("R3").Value = 3025.976557666797
$ws.Range("S3").Value = 0.06538012207349464
$ws.Range("T3").Value = 0.06538012207349464

$ws.Range("G4").Value = 14.734287
$ws.Range("H4").Value = 44.202861
$ws.Range("I4").Value = 0.4000023944294819
$ws.Range("J4").Value = 0.400002394429482
$ws.Range("M4").Value = 58.02175166666666
$ws.Range("N4").Value = 174.065255
$ws.Range("O4").Value = 0.4156043142904646
$ws.Range("P4").Value = 0.4156043142904646
$ws.Range("Q4").Value = 854.9091412993949
$ws.Range("R4").Value = 7694.182271694553
$ws.Range("S4").Value = 0.1662427208514088
$ws.Range("T4").Value = 0.1662427208514088

$ws.Range("G5").Value = 14.734287
$ws.Range("H5").Value = 44.202861
$ws.Range("I5").Value = 0.4000023944294819
$ws.Range("J5").Value = 0.400002394429482
$ws.Range("M5").Value = 15.16934033333333
$ws.Range("N5").Value = 45.508021
$ws.Range("O5").Value = 0.1086565487318021
$ws.Range("P5").Value = 0.1086565487318021
$ws.Range("Q5").Value = 223.509414072009
$ws.Range("R5").Value = 2011.584726648081
$ws.Range("S5").Value = 0.04346287966316454
$ws.Range("T5").Value = 0.04346287966316455

$ws.Range("G6").Value = 14.734287
$ws.Range("H6").Value = 44.202861
$ws.Range("I6").Value = 0.4000023944294819
$ws.Range("J6").Value = 0.400002394429482
$ws.Range("M6").Value = 41.36709099999999
$ws.Range("N6").Value = 124.101273
$ws.Range("O6").Value = 0.2963085566257249
$ws.Range("P6").Value = 0.2963085566257249
$ws.Range("Q6").Value = 609.5145911491169
$ws.Range("R6").Value = 5485.631320342052
$ws.Range("S6").Value = 0.1185241321402337
$ws.Range("T6").Value = 0.1185241321402337

$ws.Range("I7").Value = 0.3923645715978801
$ws.Range("J7").Value = 0.3923645715978802
$ws.Range("M7").Value = 2.231113333333334
$ws.Range("N7").Value = 6.69334
$ws.Range("O7").Value = 0.01598125358798882
$ws.Range("P7").Value = 0.01598125358798882
$ws.Range("Q7").Value = 32.24615606432
$ws.Range("R7").Value = 290.21540457888
$ws.Range("S7").Value = 0.006270477717648318
$ws.Range("T7").Value = 0.006270477717648319

$ws.Range("I8").Value = 0.3923645715978801
$ws.Range("J8").Value = 0.3923645715978802
$ws.Range("O8").Value = 0.1634493267640196
$ws.Range("P8").Value = 0.1634493267640195
$ws.Range("S8").Value = 0.06413172507372647
$ws.Range("T8").Value = 0.06413172507372647

$ws.Range("I9").Value = 0.3923645715978801
$ws.Range("J9").Value = 0.3923645715978802
$ws.Range("M9").Value = 58.02175166666666
$ws.Range("N9").Value = 174.065255
$ws.Range("O9").Value = 0.4156043142904646
$ws.Range("P9").Value = 0.4156043142904646
$ws.Range("Q9").Value = 838.58512762024
$ws.Range("R9").Value = 7547.266148582159
$ws.Range("S9").Value = 0.1630684087308089
$ws.Range("T9").Value = 0.1630684087308089

$ws.Range("I10").Value = 0.3923645715978801
$ws.Range("J10").Value = 0.3923645715978802
$ws.Range("M10").Value = 15.16934033333333
$ws.Range("N10").Value = 45.508021
$ws.Range("O10").Value = 0.1086565487318021
$ws.Range("P10").Value = 0.1086565487318021
$ws.Range("Q10").Value = 219.241626354608
$ws.Range("R10").Value = 1973.174637191472
$ws.Range("S10").Value = 0.04263298019445773
$ws.Range("T10").Value = 0.04263298019445774

$ws.Range("I11").Value = 0.3923645715978801
$ws.Range("J11").Value = 0.3923645715978802
$ws.Range("M11").Value = 41.36709099999999
$ws.Range("N11").Value = 124.101273
$ws.Range("O11").Value = 0.2963085566257249
$ws.Range("P11").Value = 0.2963085566257249
$ws.Range("Q11").Value = 597.876249665904
$ws.Range("R11").Value = 5380.886246993136
$ws.Range("S11").Value = 0.1162609798812388
$ws.Range("T11").Value = 0.1162609798812388

$ws.Range("G12").Value = 1.259379333333333
$ws.Range("H12").Value = 3.778138
$ws.Range("I12").Value = 0.03418928576783783
$ws.Range("J12").Value = 0.03418928576783784
$ws.Range("M12").Value = 2.231113333333334
$ws.Range("N12").Value = 6.69334
$ws.Range("O12").Value = 0.01598125358798882
$ws.Range("P12").Value = 0.01598125358798882
$ws.Range("Q12").Value = 2.809818022324445
$ws.Range("R12").Value = 25.28836220092
$ws.Range("S12").Value = 0.0005463876458480334
$ws.Range("T12").Value = 0.0005463876458480335

$ws.Range("G13").Value = 1.259379333333333
$ws.Range("H13").Value = 3.778138
$ws.Range("I13").Value = 0.03418928576783783
$ws.Range("J13").Value = 0.03418928576783784
$ws.Range("O13").Value = 0.1634493267640196
$ws.Range("P13").Value = 0.1634493267640195
$ws.Range("Q13").Value = 28.73759943484734
$ws.Range("R13").Value = 258.638394913626
$ws.Range("S13").Value = 0.00558821574129577
$ws.Range("T13").Value = 0.00558821574129577

$ws.Range("G14").Value = 1.259379333333333
$ws.Range("H14").Value = 3.778138
$ws.Range("I14").Value = 0.03418928576783783
$ws.Range("J14").Value = 0.03418928576783784
$ws.Range("M14").Value = 58.02175166666666
$ws.Range("N14").Value = 174.065255
$ws.Range("O14").Value = 0.4156043142904646
$ws.Range("P14").Value = 0.4156043142904646
$ws.Range("Q14").Value = 73.07139493279888
$ws.Range("R14").Value = 657.6425543951899
$ws.Range("S14").Value = 0.01420921466762298
$ws.Range("T14").Value = 0.01420921466762299

$ws.Range("G15").Value = 1.259379333333333
$ws.Range("H15").Value = 3.778138
$ws.Range("I15").Value = 0.03418928576783783
$ws.Range("J15").Value = 0.03418928576783784
$ws.Range("M15").Value = 15.16934033333333
$ws.Range("N15").Value = 45.508021
$ws.Range("O15").Value = 0.1086565487318021
$ws.Range("P15").Value = 0.1086565487318021
$ws.Range("Q15").Value = 19.10395371609978
$ws.Range("R15").Value = 171.935583444898
$ws.Range("S15").Value = 0.003714889795138581
$ws.Range("T15").Value = 0.003714889795138581

$ws.Range("G16").Value = 1.259379333333333
$ws.Range("H16").Value = 3.778138
$ws.Range("I16").Value = 0.03418928576783783
$ws.Range("J16").Value = 0.03418928576783784
$ws.Range("M16").Value = 41.36709099999999
$ws.Range("N16").Value = 124.101273
$ws.Range("O16").Value = 0.2963085566257249
$ws.Range("P16").Value = 0.2963085566257249
$ws.Range("Q16").Value = 52.09685948551933
$ws.Range("R16").Value = 468.871735369674
$ws.Range("S16").Value = 0.01013057791793247
$ws.Range("T16").Value = 0.01013057791793247

$ws.Range("G17").Value = 4.524255666666667
$ws.Range("H17").Value = 13.572767
$ws.Range("I17").Value = 0.1228232556945456
$ws.Range("J17").Value = 0.1228232556945456
$ws.Range("M17").Value = 2.231113333333334
$ws.Range("N17").Value = 6.69334
$ws.Range("O17").Value = 0.01598125358798882
$ws.Range("P17").Value = 0.01598125358798882
$ws.Range("Q17").Value = 10.09412714130889
$ws.Range("R17").Value = 90.84714427178001
$ws.Range("S17").Value = 0.001962869595756925
$ws.Range("T17").Value = 0.001962869595756925

$ws.Range("G18").Value = 4.524255666666667
$ws.Range("H18").Value = 13.572767
$ws.Range("I18").Value = 0.1228232556945456
$ws.Range("J18").Value = 0.1228232556945456
$ws.Range("O18").Value = 0.1634493267640196
$ws.Range("P18").Value = 0.1634493267640195
$ws.Range("Q18").Value = 103.2383521376177
$ws.Range("R18").Value = 929.145169238559
$ws.Range("S18").Value = 0.02007537845423851
$ws.Range("T18").Value = 0.02007537845423851

$ws.Range("G19").Value = 4.524255666666667
$ws.Range("H19").Value = 13.572767
$ws.Range("I19").Value = 0.1228232556945456
$ws.Range("J19").Value = 0.1228232556945456
$ws.Range("M19").Value = 58.02175166666666
$ws.Range("N19").Value = 174.065255
$ws.Range("O19").Value = 0.4156043142904646
$ws.Range("P19").Value = 0.4156043142904646
$ws.Range("Q19").Value = 262.5052387678427
$ws.Range("R19").Value = 2362.547148910585
$ws.Range("S19").Value = 0.05104587496185401
$ws.Range("T19").Value = 0.05104587496185403

$ws.Range("G20").Value = 4.524255666666667
$ws.Range("H20").Value = 13.572767
$ws.Range("I20").Value = 0.1228232556945456
$ws.Range("J20").Value = 0.1228232556945456
$ws.Range("M20").Value = 15.16934033333333
$ws.Range("N20").Value = 45.508021
$ws.Range("O20").Value = 0.1086565487318021
$ws.Range("P20").Value = 0.1086565487318021
$ws.Range("Q20").Value = 68.62997396267856
$ws.Range("R20").Value = 617.669765664107
$ws.Range("S20").Value = 0.01334555106777298
$ws.Range("T20").Value = 0.01334555106777299

$ws.Range("G21").Value = 4.524255666666667
$ws.Range("H21").Value = 13.572767
$ws.Range("I21").Value = 0.1228232556945456
$ws.Range("J21").Value = 0.1228232556945456
$ws.Range("M21").Value = 41.36709099999999
$ws.Range("N21").Value = 124.101273
$ws.Range("O21").Value = 0.2963085566257249
$ws.Range("P21").Value = 0.2963085566257249
$ws.Range("Q21").Value = 187.1552958702657
$ws.Range("R21").Value = 1684.397662832391
$ws.Range("S21").Value = 0.03639358161492315
$ws.Range("T21").Value = 0.03639358161492316

$ws.Range("G22").Value = 1.864631
$ws.Range("H22").Value = 5.593893
$ws.Range("I22").Value = 0.05062049251025444
$ws.Range("J22").Value = 0.05062049251025445
$ws.Range("M22").Value = 2.231113333333334
$ws.Range("N22").Value = 6.69334
$ws.Range("O22").Value = 0.01598125358798882
$ws.Range("P22").Value = 0.01598125358798882
$ws.Range("Q22").Value = 4.160203085846667
$ws.Range("R22").Value = 37.44182777261999
$ws.Range("S22").Value = 0.0008089789275552648
$ws.Range("T22").Value = 0.0008089789275552649

$ws.Range("G23").Value = 1.864631
$ws.Range("H23").Value = 5.593893
$ws.Range("I23").Value = 0.05062049251025444
$ws.Range("J23").Value = 0.05062049251025445
$ws.Range("O23").Value = 0.1634493267640196
$ws.Range("P23").Value = 0.1634493267640195
$ws.Range("Q23").Value = 42.54875187602899
$ws.Range("R23").Value = 382.9387668842609
$ws.Range("S23").Value = 0.008273885421264183
$ws.Range("T23").Value = 0.008273885421264183

$ws.Range("G24").Value = 1.864631
$ws.Range("H24").Value = 5.593893
$ws.Range("I24").Value = 0.05062049251025444
$ws.Range("J24").Value = 0.05062049251025445
$ws.Range("M24").Value = 58.02175166666666
$ws.Range("N24").Value = 174.065255
$ws.Range("O24").Value = 0.4156043142904646
$ws.Range("P24").Value = 0.4156043142904646
$ws.Range("Q24").Value = 108.1891568319683
$ws.Range("R24").Value = 973.7024114877148
$ws.Range("S24").Value = 0.02103809507876989
$ws.Range("T24").Value = 0.0210380950787699

$ws.Range("G25").Value = 1.864631
$ws.Range("H25").Value = 5.593893
$ws.Range("I25").Value = 0.05062049251025444
$ws.Range("J25").Value = 0.05062049251025445
$ws.Range("M25").Value = 15.16934033333333
$ws.Range("N25").Value = 45.508021
$ws.Range("O25").Value = 0.1086565487318021
$ws.Range("P25").Value = 0.1086565487318021
$ws.Range("Q25").Value = 28.28522223508367
$ws.Range("R25").Value = 254.567000115753
$ws.Range("S25").Value = 0.005500248011268286
$ws.Range("T25").Value = 0.005500248011268287

$ws.Range("G26").Value = 1.864631
$ws.Range("H26").Value = 5.593893
$ws.Range("I26").Value = 0.05062049251025444
$ws.Range("J26").Value = 0.05062049251025445
$ws.Range("M26").Value = 41.36709099999999
$ws.Range("N26").Value = 124.101273
$ws.Range("O26").Value = 0.2963085566257249
$ws.Range("P26").Value = 0.2963085566257249
$ws.Range("Q26").Value = 77.13436025842098
$ws.Range("R26").Value = 694.2092423257889
$ws.Range("S26").Value = 0.01499928507139681
$ws.Range("T26").Value = 0.01499928507139681
